$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Mark "Chapter 17: Objects" (row 18) and "Chapter 18: Classes" (row 19) as DONE
# (copy value + formatting from an existing "DONE" cell so the green fill style is reused)
$ws.Range("B2").Copy($ws.Range("B18"))
$ws.Range("B2").Copy($ws.Range("B19"))

# Update the selection to match the diff (F19)
$ws.Range("F19").Select()
